# Rerun of the analysis script: the ml_results table now also reports the
# confidence interval (ci.lower / ci.upper) for every estimated row, so two
# new columns are appended after "p.value".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# Confidence-interval bounds for each data row (rows with no SE/df/t/p in the
# original run - q_b11_b21, q_rxy1_rxy2, cross_over_point - likewise have no
# ci.lower/ci.upper and are left untouched).
$ciValues = @(
    @{ Row = 2;  Lower = 0.190999578934371;   Upper = 0.662429415014207 }
    @{ Row = 3;  Lower = -0.103850581375686;  Upper = 0.0169399409349519 }
    @{ Row = 4;  Lower = -0.185344451030458;  Upper = -0.0364881611523818 }
    @{ Row = 5;  Lower = -0.44062342489154;   Upper = 0.0718737891819489 }
    @{ Row = 6;  Lower = -0.604809899969479;  Upper = -0.11906696409819 }
    @{ Row = 7;  Lower = -0.383112909652411;  Upper = 0.0624927658079422 }
    @{ Row = 8;  Lower = -0.683750162797184;  Upper = -0.134607677701725 }
    @{ Row = 9;  Lower = -0.142350542755016;  Upper = -0.0120210835567704 }
    @{ Row = 10; Lower = 0.0176799648176695;  Upper = 0.0913996322421703 }
    @{ Row = 11; Lower = 0.0301958803538025;  Upper = 0.104726091388303 }
    @{ Row = 15; Lower = -0.0709049735742826; Upper = 0.0514553190046014 }
    @{ Row = 16; Lower = -0.261573988080058;  Upper = 0.189822692562809 }
    @{ Row = 17; Lower = -0.381091527028979;  Upper = 0.189905505498427 }
    @{ Row = 18; Lower = -0.207701162751372;  Upper = 0.0338798818699037 }
    @{ Row = 19; Lower = -0.766225819304821;  Upper = 0.124985531615884 }
    @{ Row = 20; Lower = -0.881246849783079;  Upper = 0.143747578363898 }
    @{ Row = 21; Lower = 0.0301958803538025;  Upper = 0.104726091388303 }
    @{ Row = 22; Lower = 0.0240421671135407;  Upper = 0.284701085510033 }
    @{ Row = 23; Lower = 0.111394962152516;   Upper = 0.386342734501924 }
    @{ Row = 24; Lower = 0.0886934331537079;  Upper = 1.05028455118967 }
    @{ Row = 25; Lower = 0.0484157167891871;  Upper = 0.306711511568891 }
    @{ Row = 26; Lower = 0.0640029561539502;  Upper = 1.02862354362331 }
)

foreach ($rowInfo in $ciValues) {
    $ws.Range("G" + $rowInfo.Row).Value = $rowInfo.Lower
    $ws.Range("H" + $rowInfo.Row).Value = $rowInfo.Upper
}
